$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $value
    $cell.ClearFormats()
}

Set-TextValue 2 4 "42.736.25"
Set-TextValue 2 5 "  -0.21%  "

Set-TextValue 3 4 "2.528.33"
Set-TextValue 3 5 "  -1.71%  "

Set-TextValue 4 4 "1.00"
Set-TextValue 4 5 "  -0.09%  "

Set-TextValue 5 4 "309.77"
Set-TextValue 5 5 "  -0.82%  "

Set-TextValue 6 4 "100.17"
Set-TextValue 6 5 "  +1.88%  "

Set-TextValue 7 5 "  -1.04%  "

Set-TextValue 8 5 "  +0.00%  "

Set-TextValue 9 5 "  -2.00%  "

Set-TextValue 10 4 "35.48"
Set-TextValue 10 5 "  -0.60%  "

Set-TextValue 11 5 "  -0.59%  "

Set-TextValue 12 4 "7.31"
Set-TextValue 12 5 "  -1.67%  "

Set-TextValue 13 5 "  +1.04%  "

Set-TextValue 14 4 "2.916.29"
Set-TextValue 14 5 "  -1.72%  "

Set-TextValue 15 2 "WrappedEther"
Set-TextValue 15 3 "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue 15 4 "2.576.19"
Set-TextValue 15 5 "  -1.01%  "

Set-TextValue 16 2 "Chainlink"
Set-TextValue 16 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue 16 4 "15.34"
Set-TextValue 16 5 "  -3.65%  "

Set-TextValue 17 4 "0.813"
Set-TextValue 17 5 "  -3.65%  "

Set-TextValue 18 4 "42.719.22"
Set-TextValue 18 5 "  -0.35%  "

Set-TextValue 19 5 "  -0.96%  "

Set-TextValue 20 4 "0.0₃0952"
Set-TextValue 20 5 "  -0.86%  "

Set-TextValue 21 4 "12.26"
Set-TextValue 21 5 "  -1.70%  "

Set-TextValue 22 4 "69.30"
Set-TextValue 22 5 "  -0.72%  "

Set-TextValue 23 4 "243.04"
Set-TextValue 23 5 "  -2.42%  "

Set-TextValue 24 5 "  -2.52%  "

Set-TextValue 25 5 "  -1.65%  "

Set-TextValue 26 5 "  -0.04%  "

Set-TextValue 27 4 "25.39"
Set-TextValue 27 5 "  -6.36%  "

Set-TextValue 28 4 "2.34"
Set-TextValue 28 5 "  -2.23%  "

Set-TextValue 29 4 "10.16"
Set-TextValue 29 5 "  -0.60%  "

Set-TextValue 30 4 "38.47"
Set-TextValue 30 5 "  -2.87%  "

Set-TextValue 31 4 "160.73"
Set-TextValue 31 5 "  +0.95%  "

Set-TextValue 32 4 "5.76"
Set-TextValue 32 5 "  -0.26%  "

Set-TextValue 33 4 "2.80"
Set-TextValue 33 5 "  +8.79%  "

Set-TextValue 34 4 "2.68"
Set-TextValue 34 5 "  +0.13%  "

Set-TextValue 35 5 "  -1.52%  "

Set-TextValue 36 4 "18.41"
Set-TextValue 36 5 "  -0.74%  "

Set-TextValue 37 5 "  -6.10%  "

Set-TextValue 38 5 "  -6.84%  "

Set-TextValue 39 5 "  -1.02%  "

Set-TextValue 40 5 "  -0.48%  "

Set-TextValue 41 2 "EnergySwap"
Set-TextValue 41 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue 41 4 "22.51"
Set-TextValue 41 5 "  -1.46%  "

Set-TextValue 42 2 "RenderToken"
Set-TextValue 42 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 42 4 "4.19"
Set-TextValue 42 5 "  +1.93%  "

Set-TextValue 43 5 "  +0.17%  "

Set-TextValue 44 2 "VeChain"
Set-TextValue 44 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 44 4 "0.0300"
Set-TextValue 44 5 "  -0.40%  "

Set-TextValue 45 2 "NEARProtocol"
Set-TextValue 45 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue 45 4 "3.28"
Set-TextValue 45 5 "  +2.02%  "

Set-TextValue 46 4 "1.996.54"
Set-TextValue 46 5 "  -0.07%  "

Set-TextValue 47 4 "8.82"
Set-TextValue 47 5 "  -2.22%  "

Set-TextValue 48 4 "2.768.49"
Set-TextValue 48 5 "  -1.78%  "

Set-TextValue 49 5 "  -3.13%  "

Set-TextValue 50 4 "79.22"
Set-TextValue 50 5 "  -2.84%  "

Set-TextValue 51 4 "71.82"
Set-TextValue 51 5 "  -3.39%  "
